# Refresh the crypto price ticker: latest Price/Volume(1h) figures
# pulled in by the scheduled GitHub Actions job, plus a reshuffle of
# the WrappedEther / Polygon / Polkadot rows (12-14) to match the new
# coinranking.com ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.481.87'
$ws.Range('E2').Value = '  +0.76%  '

$ws.Range('D3').Value = '1.880.92'
$ws.Range('E3').Value = '  +1.23%  '

$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = "'0.7171"
$ws.Range('E5').Value = '  +1.44%  '

$ws.Range('D6').Value = "'242.17"
$ws.Range('E6').Value = '  +1.70%  '

$ws.Range('D7').Value = "'1.000"
$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = "'0.07866"
$ws.Range('E8').Value = '  -1.72%  '

$ws.Range('D9').Value = "'0.3130"
$ws.Range('E9').Value = '  +3.45%  '

$ws.Range('D10').Value = "'25.18"

$ws.Range('D11').Value = "'0.08263"
$ws.Range('E11').Value = '  +1.00%  '

$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = "'0.7328"
$ws.Range('E12').Value = '  +3.77%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.294"
$ws.Range('E13').Value = '  +1.93%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.862.43'
$ws.Range('E14').Value = '  +0.03%  '

$ws.Range('D15').Value = "'91.30"
$ws.Range('E15').Value = '  +1.77%  '

$ws.Range('D16').Value = '29.512.50'
$ws.Range('E16').Value = '  +0.83%  '

$ws.Range('D17').Value = "'5.940"
$ws.Range('E17').Value = '  +2.26%  '

$ws.Range('D18').Value = "'248.45"
$ws.Range('E18').Value = '  +4.06%  '

$ws.Range('D19').Value = "'0.000007896"
$ws.Range('E19').Value = '  -0.50%  '

$ws.Range('D20').Value = "'13.31"
$ws.Range('E20').Value = '  +0.50%  '

$ws.Range('D21').Value = "'0.9997"
$ws.Range('E21').Value = '  +0.13%  '

$ws.Range('D22').Value = "'7.979"
$ws.Range('E22').Value = '  +6.69%  '

$ws.Range('D23').Value = "'1.000"
$ws.Range('E23').Value = '  +0.06%  '

$ws.Range('D24').Value = "'0.1596"
$ws.Range('E24').Value = '  +11.32%  '

$ws.Range('D25').Value = "'164.19"
$ws.Range('E25').Value = '  +0.78%  '

$ws.Range('D26').Value = "'9.048"
$ws.Range('E26').Value = '  +1.86%  '

$ws.Range('D27').Value = "'18.34"
$ws.Range('E27').Value = '  +1.30%  '

$ws.Range('D28').Value = "'1.362"
$ws.Range('E28').Value = '  -4.64%  '

$ws.Range('D29').Value = "'1.493"
$ws.Range('E29').Value = '  +1.21%  '

$ws.Range('D30').Value = "'4.380"
$ws.Range('E30').Value = '  +0.11%  '

$ws.Range('D31').Value = "'4.136"
$ws.Range('E31').Value = '  +2.78%  '

$ws.Range('D32').Value = "'0.05321"
$ws.Range('E32').Value = '  +2.46%  '

$ws.Range('D33').Value = "'1.940"
$ws.Range('E33').Value = '  +0.66%  '

$ws.Range('E34').Value = '  +3.70%  '

$ws.Range('D35').Value = "'0.7244"
$ws.Range('E35').Value = '  +1.40%  '

$ws.Range('D36').Value = "'2.678"
$ws.Range('E36').Value = '  +1.04%  '

$ws.Range('D37').Value = "'0.01872"
$ws.Range('E37').Value = '  +0.95%  '

$ws.Range('D38').Value = '1.265.92'
$ws.Range('E38').Value = '  +11.14%  '

$ws.Range('D39').Value = "'2.727"
$ws.Range('E39').Value = '  +0.15%  '

$ws.Range('D40').Value = "'0.9118"
$ws.Range('E40').Value = '  -2.60%  '

$ws.Range('D41').Value = "'74.30"
$ws.Range('E41').Value = '  +5.42%  '

$ws.Range('D42').Value = "'6.110"
$ws.Range('E42').Value = '  +2.62%  '

$ws.Range('E43').Value = '  +0.07%  '

$ws.Range('D44').Value = "'103.77"
$ws.Range('E44').Value = '  +1.01%  '

$ws.Range('D45').Value = '2.030.56'
$ws.Range('E45').Value = '  +1.05%  '

$ws.Range('D46').Value = "'0.5327"

$ws.Range('E47').Value = '  +0.44%  '

$ws.Range('D48').Value = "'2.920"
$ws.Range('E48').Value = '  +13.23%  '

$ws.Range('E49').Value = '  -0.08%  '

$ws.Range('D50').Value = "'0.4337"
$ws.Range('E50').Value = '  +1.63%  '

$ws.Range('D51').Value = "'9.292"
$ws.Range('E51').Value = '  +1.22%  '
